# Auto-generated: apply scheduled-runner market data refresh to Hyperion_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7259.483
$ws.Range("I74").Value = 6288.75
$ws.Range("J74").Value = 7629.2856
$ws.Range("K74").Value = 6288.75
$ws.Range("L74").Value = 7629.2856
$ws.Range("M74").Value = -5352.75
$ws.Range("N74").Value = -9501.2856

$ws.Range("H77").Value = 7259.483
$ws.Range("I77").Value = 6288.75
$ws.Range("J77").Value = 7629.2856
$ws.Range("K77").Value = 31443.75
$ws.Range("L77").Value = 38146.428
$ws.Range("M77").Value = -26763.75
$ws.Range("N77").Value = -47506.428

$ws.Range("H106").Value = 33335360
$ws.Range("I106").Value = 35716316
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 35716316
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -35715685
$ws.Range("N106").Value = -3262

$ws.Range("H107").Value = 22289134
$ws.Range("I107").Value = 23809786
$ws.Range("K107").Value = 23809786
$ws.Range("M107").Value = -23807866

$ws.Range("H141").Value = 54250
$ws.Range("I141").Value = 26800
$ws.Range("K141").Value = 80400
$ws.Range("M141").Value = -75220


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2922.0417
$ws.Range("I132").Value = 1804.4
$ws.Range("J132").Value = 4784.778
$ws.Range("K132").Value = 5413.200000000001
$ws.Range("L132").Value = 14354.334
$ws.Range("M132").Value = -2883.200000000001
$ws.Range("N132").Value = -19414.334

$ws.Range("H140").Value = 82666
$ws.Range("J140").Value = 82666
$ws.Range("L140").Value = 82666
$ws.Range("N140").Value = -93026


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3451290.8
$ws.Range("I86").Value = 5002926
$ws.Range("J86").Value = 3212.4443
$ws.Range("K86").Value = 5002926
$ws.Range("L86").Value = 3212.4443
$ws.Range("M86").Value = -5001803
$ws.Range("N86").Value = -5458.4443

$ws.Range("H89").Value = 3451290.8
$ws.Range("I89").Value = 5002926
$ws.Range("J89").Value = 3212.4443
$ws.Range("K89").Value = 25014630
$ws.Range("L89").Value = 16062.2215
$ws.Range("M89").Value = -25009014
$ws.Range("N89").Value = -27294.2215

$ws.Range("H138").Value = 79846.336
$ws.Range("J138").Value = 79846.336
$ws.Range("L138").Value = 79846.336
$ws.Range("N138").Value = -90126.336


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 835.5
$ws.Range("I16").Value = 692.4211
$ws.Range("J16").Value = 1741.6666
$ws.Range("K16").Value = 692.4211
$ws.Range("L16").Value = 1741.6666
$ws.Range("M16").Value = -405.4211
$ws.Range("N16").Value = -2315.6666

$ws.Range("H31").Value = 21443.547
$ws.Range("I31").Value = 3687.2144
$ws.Range("J31").Value = 27817.615
$ws.Range("K31").Value = 3687.2144
$ws.Range("L31").Value = 27817.615
$ws.Range("M31").Value = -3392.2144
$ws.Range("N31").Value = -28407.615

$ws.Range("H34").Value = 21443.547
$ws.Range("I34").Value = 3687.2144
$ws.Range("J34").Value = 27817.615
$ws.Range("K34").Value = 3687.2144
$ws.Range("L34").Value = 27817.615
$ws.Range("M34").Value = -3485.2144
$ws.Range("N34").Value = -28221.615

$ws.Range("H113").Value = 835.5
$ws.Range("I113").Value = 692.4211
$ws.Range("J113").Value = 1741.6666
$ws.Range("K113").Value = 692.4211
$ws.Range("L113").Value = 1741.6666
$ws.Range("M113").Value = 1477.5789
$ws.Range("N113").Value = -6081.6666

$ws.Range("H134").Value = 4663.3335

$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

$ws.Range("H132").Value = 3012.8235
$ws.Range("I132").Value = 2812.2964
$ws.Range("K132").Value = 8436.889200000001
$ws.Range("M132").Value = -5906.889200000001

$ws.Range("H133").Value = 59994.332
$ws.Range("J133").Value = 59994.332
$ws.Range("L133").Value = 59994.332
$ws.Range("N133").Value = -70114.332

$ws.Range("H136").Value = 12959.517
$ws.Range("J136").Value = 12959.517
$ws.Range("L136").Value = 38878.551
$ws.Range("N136").Value = -43978.551

$ws.Range("H140").Value = 82498.8
$ws.Range("J140").Value = 82498.8
$ws.Range("L140").Value = 82498.8
$ws.Range("N140").Value = -92858.8

$ws.Range("H141").Value = 54122.57
$ws.Range("J141").Value = 65771.6
$ws.Range("L141").Value = 65771.6
$ws.Range("N141").Value = -76131.6


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4392.6597
$ws.Range("I7").Value = 2964.7646
$ws.Range("K7").Value = 2964.7646
$ws.Range("M7").Value = -2852.7646

$ws.Range("H17").Value = 19932.445
$ws.Range("I17").Value = 17199.143
$ws.Range("K17").Value = 17199.143
$ws.Range("M17").Value = -17029.143

$ws.Range("H36").Value = 89715
$ws.Range("J36").Value = 89715
$ws.Range("L36").Value = 89715
$ws.Range("N36").Value = -90839

$ws.Range("H126").Value = 4392.6597
$ws.Range("I126").Value = 2964.7646
$ws.Range("K126").Value = 8894.2938
$ws.Range("M126").Value = -6424.293799999999

$ws.Range("H132").Value = 8276.575
$ws.Range("I132").Value = 8526.516
$ws.Range("K132").Value = 25579.548
$ws.Range("M132").Value = -23049.548

$ws.Range("H136").Value = 43448.348
$ws.Range("I136").Value = 98982.81
$ws.Range("J136").Value = 5828.2256
$ws.Range("K136").Value = 296948.43
$ws.Range("L136").Value = 17484.6768
$ws.Range("M136").Value = -294398.43
$ws.Range("N136").Value = -22584.6768

$ws.Range("H139").Value = 75799.625
$ws.Range("J139").Value = 75799.625
$ws.Range("L139").Value = 75799.625
$ws.Range("N139").Value = -86079.625

$ws.Range("H140").Value = 102757
$ws.Range("J140").Value = 128565
$ws.Range("L140").Value = 128565
$ws.Range("N140").Value = -138925


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 27786506
$ws.Range("J81").Value = 10996.667
$ws.Range("L81").Value = 21993.334
$ws.Range("N81").Value = -24115.334

$ws.Range("H84").Value = 27786506
$ws.Range("J84").Value = 10996.667
$ws.Range("L84").Value = 109966.67
$ws.Range("N84").Value = -120574.67

$ws.Range("H92").Value = 57500
$ws.Range("J92").Value = 57500
$ws.Range("L92").Value = 57500
$ws.Range("N92").Value = -62492

$ws.Range("H132").Value = 15319683
$ws.Range("I132").Value = 20411866
$ws.Range("K132").Value = 61235598
$ws.Range("M132").Value = -61233068

$ws.Range("H141").Value = 99178.75
$ws.Range("J141").Value = 99178.75
$ws.Range("L141").Value = 99178.75
$ws.Range("N141").Value = -109538.75

